$wb = $excel.ActiveWorkbook

# Rename the second worksheet from "sadgdda" to "sgfh"
$ws = $wb.Worksheets.Item("sadgdda")
$ws.Name = "sgfh"

# Updated data grid (recomputed after changing the slope comparison threshold to 1/3
# and the RT interval to 380 ms)
$csv = @"
100,2273,2273,2272,0,0,100,100,0,34.12737035751343
101,1864,1865,1863,0,1,99.94635193133047,100,0.0005361930294906167,43.32527709007263
103,2083,2084,2082,0,1,99.951992318771,100,0.0004798464491362764,39.4021680355072
105,2595,2572,2566,28,5,99.80552314274601,98.92058596761758,0.01283048211508554,35.23060369491577
106,2026,2027,2024,1,2,99.90128331688055,99.95061728395062,0.001480019733596448,40.48703908920288
108,1761,1763,1755,5,7,99.60272417707151,99.71590909090909,0.006806579693703914,34.69096851348877
109,2533,2532,2531,1,0,100,99.96050552922591,0.0003949447077409163,38.86244583129883
111,2124,2124,2122,1,1,99.95289684408856,99.95289684408856,0.0009416195856873823,33.89425992965698
112,2540,2539,2538,1,0,100,99.96061441512407,0.0003938558487593541,37.50945687294006
113,1794,1795,1793,0,1,99.94425863991081,100,0.0005571030640668524,34.1624219417572
114,1861,1879,1859,1,19,98.98828541001065,99.94623655913979,0.0106439595529537,29.7993586063385
115,1953,1953,1952,0,0,100,100,0,33.37823796272278
116,2389,2412,2388,0,23,99.0460389879718,100,0.009535655058043118,31.94921398162842
117,1536,1535,1534,1,0,100,99.93485342019544,0.0006514657980456026,32.20031881332397
118,2284,2278,2277,6,0,100,99.73718791064388,0.002633889376646181,33.36748504638672
119,1994,1987,1986,7,0,100,99.64877069744104,0.003522898842476095,36.16890358924866
121,1861,1863,1860,0,2,99.89258861439312,100,0.001073537305421363,37.6081337928772
122,2476,2476,2475,0,0,100,100,0,35.98488140106201
123,1518,1518,1517,0,0,100,100,0,36.45637249946594
124,1613,1619,1612,0,6,99.62917181705809,100,0.003705991352686844,29.39400672912598
200,2602,2601,2597,4,3,99.88461538461539,99.84621299500192,0.002691272587466359,34.79302382469177
201,1957,1963,1955,1,7,99.64322120285424,99.94887525562372,0.004075394803871625,40.05390501022339
202,2134,2136,2133,0,2,99.90632318501171,100,0.0009363295880149813,30.98800730705261
203,2962,2980,2946,15,33,98.89224572004028,99.49341438703141,0.01610738255033557,36.10018086433411
205,2646,2656,2645,0,10,99.62335216572505,100,0.003765060240963855,34.98470568656921
207,1852,1860,1847,4,12,99.35449166218397,99.78390059427336,0.008602150537634409,37.47321534156799
208,2945,2955,2942,2,12,99.5937711577522,99.9320652173913,0.004737732656514382,35.18442392349243
209,2988,3005,2986,1,18,99.40079893475367,99.96652159357214,0.00632279534109817,33.44931840896606
210,2637,2650,2631,5,18,99.32049830124575,99.8103186646434,0.008679245283018867,35.60828542709351
212,2748,2748,2747,0,0,100,100,0,38.20660948753357
213,3249,3251,3248,0,2,99.93846153846154,100,0.0006151953245155337,35.74511790275574
214,2263,2262,2258,4,3,99.86731534719151,99.82316534040672,0.003094606542882405,37.57099914550781
215,3362,3363,3361,0,1,99.97025580011898,100,0.0002973535533749628,33.86595511436462
219,2154,2154,2153,0,0,100,100,0,37.29110455513
220,2048,2048,2047,0,0,100,100,0,42.92071914672852
221,2426,2427,2425,0,1,99.95877988458368,100,0.0004120313143798929,34.92539954185486
222,2478,2483,2476,1,6,99.75825946817083,99.95962858296326,0.002819170358437374,35.46833920478821
223,2605,2605,2604,0,0,100,100,0,30.95150446891785
228,2071,2053,2051,19,1,99.95126705653021,99.08212560386474,0.00974184120798831,38.36787414550781
230,2257,2256,2255,1,0,100,99.95567375886525,0.0004432624113475177,39.30532479286194
231,1571,1571,1570,0,0,100,100,0,36.05975842475891
232,1780,1780,1778,1,1,99.94378864530636,99.94378864530636,0.001123595505617978,35.33758568763733
233,3078,3079,3077,0,1,99.96751137102014,100,0.0003247807729782397,40.8732316493988
234,2753,2753,2752,0,0,100,100,0,35.62307405471802
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$rowCount = $lines.Count
$colCount = 10

$data = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    $fields = $lines[$r].Trim() -split ","
    for ($c = 0; $c -lt $colCount; $c++) {
        $data[$r,$c] = [double]$fields[$c]
    }
}

$rng = $ws.Range("A1:J44")
$rng.Value2 = $data
